$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update Runmode column: all rows now run as "Y"
$ws.Range("C2:C7").Value = "Y"

# Update the selection to reflect the active range used when making the edit
$ws.Range("C2:C7").Select()
